$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 13.191
$ws.Range("C12").Value = -13.052
$ws.Range("E12").Value = 13.086
$ws.Range("E14").Value = 13.072
$ws.Range("E22").Value = 13.126
